$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column values are written as exact text, preserving
# trailing zeros / thousands-dot formatting instead of being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.751.88'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.603.79'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.01'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0620'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.248'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.74'
$ws.Range("E10").Value = '  +1.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.829.08'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.601.19'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  +1.18%  '
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.08'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0743'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '210.33'
$ws.Range("E18").Value = '  +0.56%  '
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.24'
$ws.Range("E22").Value = '  -4.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.07'
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.74'
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.12'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.37'
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  +0.98%  '
$ws.Range("E32").Value = '  +0.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.292.91'
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("E34").Value = '  +1.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.22'
$ws.Range("E35").Value = '  +17.38%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.592'
$ws.Range("E37").Value = '  -4.01%  '
$ws.Range("E38").Value = '  -0.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.832'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("E41").Value = '  -0.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.780'
$ws.Range("E42").Value = '  -0.48%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '63.23'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.740.98'
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.41'
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.45'
$ws.Range("E46").Value = '  +25.79%  '
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  -2.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.103'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0514'
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.57'
$ws.Range("E51").Value = '  +2.85%  '
